# Update latest output (run 194)
# Applies the optimisation_result.xlsx refresh: new "Schedule" pump window
# (row 6 added), revised Schedule metrics, and revised "Detailed" price /
# classification / pump-status columns.

$wb = $excel.ActiveWorkbook
$wsSchedule = $wb.Worksheets.Item(1)
$wsDetailed = $wb.Worksheets.Item(2)

# ---------------------------------------------------------------------
# Sheet "Schedule": rows 2-5 get revised values, row 6 is brand new.
# (A1:F5 -> A1:F6)
# ---------------------------------------------------------------------
$scheduleRows = @(
    @{Row=2; A=46076;              B=46076.20833333334; C=5; D=18.9;  E=774.0742035;        F=40.95630706349206},
    @{Row=3; A=46076.375;          B=46076.66666666666; C=7; D=26.46; E=695.4520169999998;  F=26.28314501133787},
    @{Row=4; A=46076.83333333334;  B=46077;              C=4; D=15.12; E=701.3367465000001; F=46.38470545634922},
    @{Row=5; A=46077.33333333334;  B=46077.66666666666; C=8; D=30.24; E=594.1234065;        F=19.64693804563492},
    @{Row=6; A=46077.83333333334;  B=46078;              C=4; D=15.12; E=660.41255475;      F=43.67807901785714}
)

foreach ($r in $scheduleRows) {
    $wsSchedule.Cells.Item($r.Row, 1).Value = $r.A
    $wsSchedule.Cells.Item($r.Row, 2).Value = $r.B
    $wsSchedule.Cells.Item($r.Row, 3).Value = $r.C
    $wsSchedule.Cells.Item($r.Row, 4).Value = $r.D
    $wsSchedule.Cells.Item($r.Row, 5).Value = $r.E
    $wsSchedule.Cells.Item($r.Row, 6).Value = $r.F
}

# New row 6's Start/Stop Time columns (A & B) need the same date-time
# display format as the rest of the column.
$wsSchedule.Range("A6:B6").NumberFormat = $wsSchedule.Range("A5:B5").NumberFormat

# ---------------------------------------------------------------------
# Sheet "Detailed": scattered Price (B) / Type (C) / Pump_Status (E)
# updates across rows 12, 19, 38-97.
# ---------------------------------------------------------------------
$detailedRows = @(
    @{Row=12; E="OFF"},
    @{Row=19; E="OFF"},
    @{Row=38; B=70.92874},
    @{Row=39; B=88.27029},
    @{Row=40; B=85.13891;            C="historical"},
    @{Row=41; B=85.66633;            C="historical"},
    @{Row=42; B=99.35961;            C="historical"; E="ON"},
    @{Row=43; B=93.23457000000001;   C="historical"; E="ON"},
    @{Row=44; C="historical"},
    @{Row=45; B=78;                  C="historical"},
    @{Row=46; B=68.61024;            C="historical"},
    @{Row=47; B=103.84773;           C="historical"},
    @{Row=48; B=103.79417;           C="historical"},
    @{Row=49; B=94.47342;            C="historical"},
    @{Row=50; B=97.72414999999999;   E="OFF"},
    @{Row=51; B=92.20686000000001;   E="OFF"},
    @{Row=52; B=93.19385;            E="OFF"},
    @{Row=53; B=93.9629;             E="OFF"},
    @{Row=54; B=92.51766000000001;   E="OFF"},
    @{Row=55; B=102.77994;           E="OFF"},
    @{Row=56; B=103.67301;           E="OFF"},
    @{Row=57; B=105},
    @{Row=58; B=105.79},
    @{Row=59; B=84.79000000000001},
    @{Row=60; B=84.79000000000001},
    @{Row=61; B=108.01},
    @{Row=62; B=108.01},
    @{Row=64; B=93.30685},
    @{Row=65; B=76.29461999999999;   E="OFF"},
    @{Row=66; B=62.22123},
    @{Row=67; B=56.98},
    @{Row=69; B=22.07},
    @{Row=71; B=34.01},
    @{Row=72; B=28.92581},
    @{Row=73; B=34.01},
    @{Row=74; B=34.01},
    @{Row=75; B=34.01},
    @{Row=76; B=37.89},
    @{Row=77; B=37.89},
    @{Row=78; B=37.89},
    @{Row=79; B=37.8903},
    @{Row=80; B=37.89},
    @{Row=81; B=37.89},
    @{Row=82; B=56.98003},
    @{Row=83; B=37.89019},
    @{Row=84; B=66.19919},
    @{Row=85; B=76.28172000000001},
    @{Row=86; B=95.73061},
    @{Row=87; B=108.89},
    @{Row=88; B=108.89},
    @{Row=90; B=108.89;              E="ON"},
    @{Row=91; E="ON"},
    @{Row=92; B=98.93407999999999;   E="ON"},
    @{Row=93; B=84.79000000000001;   E="ON"},
    @{Row=94; B=70.57161000000001;   E="ON"},
    @{Row=95; B=71.47271000000001;   E="ON"},
    @{Row=96; B=72.79781;            E="ON"},
    @{Row=97; B=64.89;               E="ON"}
)

foreach ($r in $detailedRows) {
    if ($r.ContainsKey("B")) {
        $wsDetailed.Cells.Item($r.Row, 2).Value = $r.B
    }
    if ($r.ContainsKey("C")) {
        $wsDetailed.Cells.Item($r.Row, 3).Value = $r.C
    }
    if ($r.ContainsKey("E")) {
        $wsDetailed.Cells.Item($r.Row, 5).Value = $r.E
    }
}
